$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value  = "Associative"
$ws.Range("F4").Value  = "Associative"
$ws.Range("F7").Value  = "Controlled, high effort, logical, rule-based"
$ws.Range("F17").Value = "Rapid, associative"
$ws.Range("F21").Value = "Associative, holistic, analytic"
$ws.Range("F23").Value = "Associative, low effort,`npragmatic"
$ws.Range("F26").Value = "Associative, contextualized, low effort, rapid"
$ws.Range("F24").Value = "Controlled, high-effort, rule-based"

$ws.Range("H24").Select() | Out-Null
